# Update cryptocurrency price/volume data per the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.108.64"
$ws.Range("E2").Value = "  +0.36%  "

# Row 3
$ws.Range("D3").Value = "3.744.99"
$ws.Range("E3").Value = "  +0.32%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'601.75"
$ws.Range("E5").Value = "  +0.18%  "

# Row 6
$ws.Range("D6").Value = "'167.46"
$ws.Range("E6").Value = "  +0.03%  "

# Row 7
$ws.Range("D7").Value = "3.743.74"
$ws.Range("E7").Value = "  +0.31%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").Value = "'0.542"
$ws.Range("E9").Value = "  +1.85%  "

# Row 10
$ws.Range("E10").Value = "  +4.06%  "

# Row 11
$ws.Range("E11").Value = "  +0.47%  "

# Row 12
$ws.Range("E12").Value = "  +0.36%  "

# Row 13
$ws.Range("D13").Value = "'38.04"
$ws.Range("E13").Value = "  +0.10%  "

# Row 14
$ws.Range("D14").Value = "'0.0000248"
$ws.Range("E14").Value = "  +2.23%  "

# Row 15
$ws.Range("D15").Value = "4.367.29"
$ws.Range("E15").Value = "  +0.24%  "

# Row 16
$ws.Range("D16").Value = "3.742.31"
$ws.Range("E16").Value = "  +0.19%  "

# Row 17
$ws.Range("D17").Value = "69.072.96"
$ws.Range("E17").Value = "  +0.32%  "

# Row 18
$ws.Range("E18").Value = "  +1.53%  "

# Row 19
$ws.Range("E19").Value = "  +0.64%  "

# Row 20
$ws.Range("E20").Value = "  -1.52%  "

# Row 21
$ws.Range("D21").Value = "'11.12"
$ws.Range("E21").Value = "  +9.15%  "

# Row 22
$ws.Range("D22").Value = "'492.65"
$ws.Range("E22").Value = "  -0.86%  "

# Row 23
$ws.Range("D23").Value = "'0.728"
$ws.Range("E23").Value = "  +0.83%  "

# Row 24
$ws.Range("E24").Value = "  +8.82%  "

# Row 25
$ws.Range("D25").Value = "'84.98"
$ws.Range("E25").Value = "  -0.17%  "

# Row 26
$ws.Range("D26").Value = "'2.29"
$ws.Range("E26").Value = "  -0.34%  "

# Row 27
$ws.Range("E27").Value = "  -0.76%  "

# Row 28
$ws.Range("E28").Value = "  -0.32%  "

# Row 29
$ws.Range("E29").Value = "  +0.02%  "

# Row 30
$ws.Range("E30").Value = "  +1.28%  "

# Row 31
$ws.Range("D31").Value = "'8.14"
$ws.Range("E31").Value = "  +2.31%  "

# Row 32
$ws.Range("E32").Value = "  +0.74%  "

# Row 33
$ws.Range("D33").Value = "'31.51"
$ws.Range("E33").Value = "  -0.55%  "

# Row 34
$ws.Range("D34").Value = "3.890.11"
$ws.Range("E34").Value = "  +0.20%  "

# Row 35
$ws.Range("D35").Value = "3.677.19"
$ws.Range("E35").Value = "  +0.36%  "

# Row 36
$ws.Range("E36").Value = "  +0.21%  "

# Row 37
$ws.Range("E37").Value = "  -0.09%  "

# Row 38
$ws.Range("D38").Value = "'1.02"
$ws.Range("E38").Value = "  +0.42%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.139"
$ws.Range("E39").Value = "  +4.99%  "

# Row 40
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'5.94"
$ws.Range("E40").Value = "  +2.28%  "

# Row 41
$ws.Range("E41").Value = "  +0.05%  "

# Row 42
$ws.Range("E42").Value = "  +6.80%  "

# Row 43
$ws.Range("D43").Value = "'48.79"
$ws.Range("E43").Value = "  -0.19%  "

# Row 44
$ws.Range("E44").Value = "  +0.32%  "

# Row 45
$ws.Range("D45").Value = "'424.60"
$ws.Range("E45").Value = "  -2.34%  "

# Row 46
$ws.Range("D46").Value = "'8.46"
$ws.Range("E46").Value = "  +0.74%  "

# Row 47
$ws.Range("E47").Value = "  +0.00%  "

# Row 48
$ws.Range("D48").Value = "'40.10"
$ws.Range("E48").Value = "  -1.08%  "

# Row 49
$ws.Range("D49").Value = "'141.97"
$ws.Range("E49").Value = "  -0.12%  "

# Row 50
$ws.Range("D50").Value = "2.781.69"
$ws.Range("E50").Value = "  +1.53%  "

# Row 51
$ws.Range("E51").Value = "  +0.37%  "
